$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-6 (by column D, L, M, N, O, P, Q, R, S, T) are cyclically
# rotated: new row2 = old row4, new row3 = old row5, new row4 = old row6,
# new row5 = old row2, new row6 = old row3.
# Capture the original values first, then write them back in rotated order.

$rows = 2..6
$data = @{}
foreach ($r in $rows) {
    $data[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

# mapping: new row -> source row
$map = @{ 2 = 4; 3 = 5; 4 = 6; 5 = 2; 6 = 3 }

foreach ($newRow in $rows) {
    $src = $data[$map[$newRow]]
    $ws.Cells.Item($newRow, 4).Value = $src.D
    $ws.Cells.Item($newRow, 12).Value = $src.L
    $ws.Cells.Item($newRow, 13).Value = $src.M
    $ws.Cells.Item($newRow, 14).Value = $src.N
    $ws.Cells.Item($newRow, 15).Value = $src.O
    $ws.Cells.Item($newRow, 16).Value = $src.P
    $ws.Cells.Item($newRow, 17).Value = $src.Q
    $ws.Cells.Item($newRow, 18).Value = $src.R
    $ws.Cells.Item($newRow, 19).Value = $src.S
    $ws.Cells.Item($newRow, 20).Value = $src.T
}
